$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for "98, 423, 780, 1105" and "98, 130, 748, 1073" swap
# places (row 2 <-> row 3), while the column headers in row 1 stay put.
# Use Copy (instead of re-typing .Value) so text-typed numeric-looking
# cells (e.g. "780", "5552") keep their original string type/format
# instead of being reinterpreted as numbers.

# Stage row 2 (A2:D2) out of the way in a scratch row far below the data.
$ws.Range("A2:D2").Copy($ws.Range("A100"))

# Move row 3's data up into row 2.
$ws.Range("A3:D3").Copy($ws.Range("A2"))

# Move the staged original row 2 data down into row 3.
$ws.Range("A100:D100").Copy($ws.Range("A3"))

# Clean up the scratch row.
$ws.Range("A100:D100").Clear()
